# Insurance Payment Avgs - add Margin % helper column (G) and rewrite
# column E (Margin) as a formula driven off it, for rows 248:281.
#
#   G<n> = (D<n> - 34) / D<n>
#   E<n> = G<n> * 100
#
# Row 248 is the "first" formula (not part of the shared group).
# Rows 249:281 are filled down as a shared formula group.
# Row 258's Margin (E258) ends up blank (its D value is blank too),
# while G258 still carries the fill-down formula (-> #DIV/0!).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E (Margin) formulas -------------------------------------------
$ws.Range("E248").Formula = "=G248*100"
$ws.Range("E249:E281").Formula = "=G249*100"

# --- Column G (new helper column) formulas ---------------------------------
$ws.Range("G248").Formula = "=(D248-34)/D248"
$ws.Range("G249:G281").Formula = "=(D249-34)/D249"

# Row 258 had its Margin formula removed by the author after the fill-down,
# leaving the cell blank while the helper column still has the formula.
$ws.Range("E258").ClearContents()

# --- Restore the view state (scroll position / active selection) ----------
$excel.ActiveWindow.ScrollRow = 213
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E258").Select()
